$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 350, shifting existing rows 350:448 down to 351:449.
$ws.Rows.Item(350).Insert()

# Populate the newly inserted row 350 with its data (same record as the
# historical row 350 but reporting a new date and a different Volumen value).
$ws.Range("A350").Value = 9
$ws.Range("B350").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C350").Value = "Metropolitana"
$ws.Range("D350").Value = 44508
$ws.Range("E350").Value = 13
$ws.Range("F350").Value = 100112023
$ws.Range("G350").Value = "Brócoli"
$ws.Range("H350").Value = "Sin especificar"
$ws.Range("I350").Value = "Primera"
$ws.Range("J350").Value = 2500
$ws.Range("K350").Value = 600
$ws.Range("L350").Value = 700
$ws.Range("M350").Value = 650
$ws.Range("N350").Value = "$/unidad"
$ws.Range("O350").Value = "Región Metropolitana"
$ws.Range("P350").Value = 650
$ws.Range("Q350").Value = 1
$ws.Range("R350").Value = "Hortaliza"
